$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 331.27274
$ws.Range("I2").Value = 331.27274
$ws.Range("K2").Value = 331.27274
$ws.Range("M2").Value = -218.27274
$ws.Range("H74").Value = 3398.7
$ws.Range("I74").Value = 2997.6
$ws.Range("K74").Value = 2997.6
$ws.Range("M74").Value = -2061.6
$ws.Range("H77").Value = 3398.7
$ws.Range("I77").Value = 2997.6
$ws.Range("K77").Value = 14988
$ws.Range("M77").Value = -10308
$ws.Range("H80").Value = 884.9286
$ws.Range("I80").Value = 1160
$ws.Range("J80").Value = 732.1111
$ws.Range("K80").Value = 3480
$ws.Range("L80").Value = 2196.3333
$ws.Range("M80").Value = -2482
$ws.Range("N80").Value = -4192.3333
$ws.Range("H83").Value = 884.9286
$ws.Range("I83").Value = 1160
$ws.Range("J83").Value = 732.1111
$ws.Range("K83").Value = 10440
$ws.Range("L83").Value = 6588.9999
$ws.Range("M83").Value = -5448
$ws.Range("N83").Value = -16572.9999
$ws.Range("H113").Value = 1946.9166
$ws.Range("I113").Value = 1926.125
$ws.Range("K113").Value = 1926.125
$ws.Range("M113").Value = 1327.875
$ws.Range("H129").Value = 746.8946999999999
$ws.Range("J129").Value = 912.25
$ws.Range("L129").Value = 2736.75
$ws.Range("N129").Value = -12736.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1370.6666
$ws.Range("I45").Value = 1256
$ws.Range("K45").Value = 1256
$ws.Range("M45").Value = -879
$ws.Range("H61").Value = 30304150
$ws.Range("I61").Value = 33334278
$ws.Range("K61").Value = 33334278
$ws.Range("M61").Value = -33334066
$ws.Range("H97").Value = 667.375
$ws.Range("I97").Value = 667.375
$ws.Range("K97").Value = 667.375
$ws.Range("M97").Value = -171.375
$ws.Range("H122").Value = 1893.3871
$ws.Range("I122").Value = 1879.0769
$ws.Range("J122").Value = 1967.8
$ws.Range("K122").Value = 5637.2307
$ws.Range("L122").Value = 5903.4
$ws.Range("M122").Value = -3187.2307
$ws.Range("N122").Value = -10803.4
$ws.Range("H132").Value = 2425.8103
$ws.Range("I132").Value = 1728.6154
$ws.Range("J132").Value = 3856.8948
$ws.Range("K132").Value = 5185.8462
$ws.Range("L132").Value = 11570.6844
$ws.Range("M132").Value = -2655.8462
$ws.Range("N132").Value = -16630.6844
$ws.Range("H136").Value = 30304150
$ws.Range("I136").Value = 33334278
$ws.Range("K136").Value = 100002834
$ws.Range("M136").Value = -100000284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 5000
$ws.Range("I19").Value = 5000
$ws.Range("K19").Value = 5000
$ws.Range("M19").Value = -4827
$ws.Range("H86").Value = 2915.7693
$ws.Range("I86").Value = 3050.2222
$ws.Range("K86").Value = 3050.2222
$ws.Range("M86").Value = -1927.2222
$ws.Range("H89").Value = 2915.7693
$ws.Range("I89").Value = 3050.2222
$ws.Range("K89").Value = 15251.111
$ws.Range("M89").Value = -9635.111000000001
$ws.Range("H99").Value = 37038028
$ws.Range("I99").Value = 47619940
$ws.Range("K99").Value = 47619940
$ws.Range("M99").Value = -47618442
$ws.Range("H110").Value = 35500
$ws.Range("J110").Value = 35500
$ws.Range("L110").Value = 35500
$ws.Range("N110").Value = -43680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 10633.667
$ws.Range("J109").Value = 10633.667
$ws.Range("L109").Value = 10633.667
$ws.Range("N109").Value = -12713.667
$ws.Range("H112").Value = 38427.363
$ws.Range("J112").Value = 38427.363
$ws.Range("L112").Value = 38427.363
$ws.Range("N112").Value = -41381.363

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1632.3334
$ws.Range("I122").Value = 596.5714
$ws.Range("J122").Value = 2538.625
$ws.Range("K122").Value = 5369.1426
$ws.Range("L122").Value = 22847.625
$ws.Range("M122").Value = -2919.1426
$ws.Range("N122").Value = -27747.625
$ws.Range("H131").Value = 24427608
$ws.Range("J131").Value = 50883.133
$ws.Range("L131").Value = 152649.399
$ws.Range("N131").Value = -162729.399
$ws.Range("H132").Value = 1153.6111
$ws.Range("I132").Value = 975.7778
$ws.Range("J132").Value = 1331.4445
$ws.Range("K132").Value = 8782.0002
$ws.Range("L132").Value = 11983.0005
$ws.Range("M132").Value = -6252.0002
$ws.Range("N132").Value = -17043.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 40551.285
$ws.Range("J42").Value = 40551.285
$ws.Range("L42").Value = 40551.285
$ws.Range("N42").Value = -41521.285
$ws.Range("H80").Value = 2993.2632
$ws.Range("I80").Value = 1573.8889
$ws.Range("J80").Value = 4270.7
$ws.Range("K80").Value = 1573.8889
$ws.Range("L80").Value = 4270.7
$ws.Range("M80").Value = -575.8888999999999
$ws.Range("N80").Value = -6266.7
$ws.Range("H83").Value = 2993.2632
$ws.Range("I83").Value = 1573.8889
$ws.Range("J83").Value = 4270.7
$ws.Range("K83").Value = 7869.4445
$ws.Range("L83").Value = 21353.5
$ws.Range("M83").Value = -2877.4445
$ws.Range("N83").Value = -31337.5
$ws.Range("H102").Value = 1487.9395
$ws.Range("I102").Value = 1427.3462
$ws.Range("K102").Value = 1427.3462
$ws.Range("M102").Value = 194.6538
$ws.Range("H109").Value = 8333.166999999999
$ws.Range("J109").Value = 8333.166999999999
$ws.Range("L109").Value = 8333.166999999999
$ws.Range("N109").Value = -10413.167
$ws.Range("H113").Value = 1236.7693
$ws.Range("I113").Value = 942
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 942
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 1228
$ws.Range("N113").Value = -6240
$ws.Range("H115").Value = 40551.285
$ws.Range("J115").Value = 40551.285
$ws.Range("L115").Value = 40551.285
$ws.Range("N115").Value = -42901.285
$ws.Range("H123").Value = 21700
$ws.Range("J123").Value = 21700
$ws.Range("L123").Value = 21700
$ws.Range("N123").Value = -26600
$ws.Range("H126").Value = 1971.7273
$ws.Range("I126").Value = 1683.3846
$ws.Range("J126").Value = 2388.2222
$ws.Range("K126").Value = 5050.1538
$ws.Range("L126").Value = 7164.6666
$ws.Range("M126").Value = -2580.1538
$ws.Range("N126").Value = -12104.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2811.4285
$ws.Range("J7").Value = 3203.2
$ws.Range("L7").Value = 3203.2
$ws.Range("N7").Value = -3427.2
$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10452
$ws.Range("H122").Value = 41668336
$ws.Range("I122").Value = 83334664
$ws.Range("K122").Value = 250003992
$ws.Range("M122").Value = -250001542
$ws.Range("H126").Value = 2811.4285
$ws.Range("J126").Value = 3203.2
$ws.Range("L126").Value = 9609.599999999999
$ws.Range("N126").Value = -14549.6
$ws.Range("H132").Value = 2979
$ws.Range("I132").Value = 3255.3333
$ws.Range("K132").Value = 9765.999899999999
$ws.Range("M132").Value = -7235.999899999999
$ws.Range("H134").Value = 35300
$ws.Range("J134").Value = 35300
$ws.Range("L134").Value = 35300
$ws.Range("N134").Value = -45440
$ws.Range("H136").Value = 1694.3334
$ws.Range("I136").Value = 1482.4286
$ws.Range("J136").Value = 2881
$ws.Range("K136").Value = 4447.2858
$ws.Range("L136").Value = 8643
$ws.Range("M136").Value = -1897.2858
$ws.Range("N136").Value = -13743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 15000
$ws.Range("J70").Value = 15000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15630
$ws.Range("H73").Value = 15000
$ws.Range("J73").Value = 15000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -17184
$ws.Range("H81").Value = 1493.2222
$ws.Range("I81").Value = 1304.875
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 2609.75
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -1548.75
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 1493.2222
$ws.Range("I84").Value = 1304.875
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 13048.75
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -7744.75
$ws.Range("N84").Value = -40608
$ws.Range("H112").Value = 12500
$ws.Range("J112").Value = 12500
$ws.Range("L112").Value = 12500
$ws.Range("N112").Value = -15454
$ws.Range("H113").Value = 525.55554
$ws.Range("J113").Value = 841.1667
$ws.Range("L113").Value = 2523.5001
$ws.Range("N113").Value = -6863.5001
$ws.Range("H133").Value = 27320
$ws.Range("J133").Value = 27320
$ws.Range("L133").Value = 27320
$ws.Range("N133").Value = -37440
